$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D5").Value = "implemented"
$ws.Range("E5").Value = "Rotation backtest implemented: Top-N momentum (with optional eligible DSL) + compare overlay UI + tests."
$ws.Range("F5").Value = "27/12/2025 03:26"
